# Applies cryptos.xlsx price/volume/row-order updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.283.26"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.551.73"
$ws.Range("E3").Value = "  -1.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.40"
$ws.Range("E5").Value = "  -1.50%  "

$ws.Range("E6").Value = "  -1.65%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.47"

$ws.Range("E9").Value = "  -2.04%  "

$ws.Range("E10").Value = "  -1.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.774.36"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.561.17"
$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.291.54"
$ws.Range("E14").Value = "  -0.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.62"
$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.508"
$ws.Range("E16").Value = "  -2.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.46"
$ws.Range("E17").Value = "  -2.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.71"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.30"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0674"
$ws.Range("E20").Value = "  -2.52%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("E22").Value = "  +1.14%  "

$ws.Range("E23").Value = "  -3.22%  "

$ws.Range("E24").Value = "  -5.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.62"
$ws.Range("E25").Value = "  -2.42%  "

$ws.Range("E26").Value = "  -1.72%  "

$ws.Range("E27").Value = "  -0.36%  "

$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("E29").Value = "  -3.11%  "

$ws.Range("E30").Value = "  -3.66%  "

$ws.Range("E31").Value = "  -4.31%  "

$ws.Range("E32").Value = "  -0.89%  "

$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.385.58"
$ws.Range("E34").Value = "  -0.49%  "

$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("E38").Value = "  -1.68%  "

$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("E40").Value = "  +1.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.509"
$ws.Range("E41").Value = "  -2.59%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0468"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.776"
$ws.Range("E44").Value = "  -1.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.41"
$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.74"
$ws.Range("E46").Value = "  -1.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.687.11"

$ws.Range("E48").Value = "  -6.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.46"
$ws.Range("E49").Value = "  -1.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.30"
$ws.Range("E50").Value = "  +6.76%  "

$ws.Range("E51").Value = "  -0.25%  "
